# Add a new "2022" column (K) to the 10.1.1 indicator table, mirroring the
# formatting already used for the "2021" column (J), and move the active
# selection the way the author left it after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (year header): K4 = 2022, formatted like J4 -----------------
$srcJ4 = $ws.Range("J4")
$k4 = $ws.Range("K4")
$srcJ4.Copy($k4)
$k4.Value2 = 2022

# --- Row 5 (bottom-40% growth rate): K5 = 3.9462868231169921 -----------
$srcJ5 = $ws.Range("J5")
$k5 = $ws.Range("K5")
$srcJ5.Copy($k5)
$k5.Value2 = 3.9462868231169921
# the new cell drops the vertical-center alignment that J5 carries
$k5.VerticalAlignment = -4107

# --- Row 6 (whole-population growth rate): K6 = 3.8007658934388928 -----
$srcJ6 = $ws.Range("J6")
$k6 = $ws.Range("K6")
$srcJ6.Copy($k6)
$k6.Value2 = 3.8007658934388928
# the new cell drops the vertical-center alignment that J6 carries
$k6.VerticalAlignment = -4107

# --- Leave the selection where the author left it after editing --------
$ws.Range("L5").Select()
